$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update JORGE's (account 005068961) balance from 286169.56 to 336169.56
$ws.Range("C3").Value = 336169.56

# 2) Remove the four trailing rows that were dropped from the bottom of the
#    list (ROBERIO -284.48, JOSE -5497.02, ALAN -7533.83, ANDRE -95721.38).
#    Do this before the insert below so the row numbers stay stable.
$ws.Rows("223:226").Delete()

# 3) Insert a new row for ROBERIO (account 005203562, balance 169.16) right
#    above the CARLOS row (account 004360431), shifting everything below it
#    down by one. The leading apostrophe forces the account number to stay
#    text (matching the other account-number cells in column A).
$ws.Rows(84).Insert()
$ws.Range("A84").Value = "'005203562"
$ws.Range("B84").Value = "ROBERIO"
$ws.Range("C84").Value = 169.16
